$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1 (Title): add " (Summarised)" run and move the _GoBack bookmark here ---
$p1 = $d.Paragraphs.Item(1)
$xml1 = '<w:p ' + $ns + '>' + `
  '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' + `
  '<w:r><w:t>Fine-Wines System User Manual</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> (Summarised)</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
$p1.Range.InsertXML($xml1) | Out-Null

# --- Paragraph 4 ("Registering A Business And Or Employee Or Administrator"): drop the bookmark (moved to title) ---
$p4 = $d.Paragraphs.Item(4)
$xml4 = '<w:p ' + $ns + '>' + `
  '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Registering </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">A Business </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>And</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> Or Employee Or Administrator</w:t></w:r>' + `
  '</w:p>'
$p4.Range.InsertXML($xml4) | Out-Null

# --- Insert a new paragraph after paragraph 7 ("Navigate to one of the areas...") ---
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item(8)
$p8.Style = "Normal"
$xml8 = '<w:p ' + $ns + '>' + `
  '<w:r><w:t>These options and rules apply for both the Wines and Grapes pages, as well as the other pages related to the management of data in the database.</w:t></w:r>' + `
  '</w:p>'
$p8.Range.InsertXML($xml8) | Out-Null

# --- Paragraph 11 ("Getting Statistics, Estimations And Graphs" -> add ", Reports") ---
$p11 = $d.Paragraphs.Item(11)
$xml11 = '<w:p ' + $ns + '>' + `
  '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' + `
  '<w:r><w:t>Getting Statistics, Estimations</w:t></w:r>' + `
  '<w:r><w:t>, Reports</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>And</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> Graphs</w:t></w:r>' + `
  '</w:p>'
$p11.Range.InsertXML($xml11) | Out-Null

# --- Paragraph 12 ("Navigate to the estimations page..." -> "Navigate to the Reports page...") ---
$p12 = $d.Paragraphs.Item(12)
$xml12 = '<w:p ' + $ns + '>' + `
  '<w:r><w:t xml:space="preserve">Navigate to the </w:t></w:r>' + `
  '<w:r><w:t>Reports</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> page and select the statistic that you would like to see (Alphabetically, Actual Production, Estimated Production, Percentage Produced). Use the relevant button or control to export the statistic.</w:t></w:r>' + `
  '</w:p>'
$p12.Range.InsertXML($xml12) | Out-Null

# --- Paragraph 13: reuse the trailing empty paragraph as the new "Checking And Managing The Stock" heading ---
$p13 = $d.Paragraphs.Item(13)
$p13.Style = "Heading 1"
$xml13 = '<w:p ' + $ns + '>' + `
  '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' + `
  '<w:r><w:t>Checking</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>And</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> Managing</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> The Stock</w:t></w:r>' + `
  '</w:p>'
$p13.Range.InsertXML($xml13) | Out-Null

# --- Paragraph 14 (new): body text about browsing stock ---
$p13 = $d.Paragraphs.Item(13)
$p13.Range.InsertParagraphAfter()
$p14 = $d.Paragraphs.Item(14)
$p14.Style = "Normal"
$xml14 = '<w:p ' + $ns + '>' + `
  '<w:r><w:t>Navigate to the Browse Stock page in order to manage and check the current stock.</w:t></w:r>' + `
  '</w:p>'
$p14.Range.InsertXML($xml14) | Out-Null
